# Update "想去人数" (want-to-go count) figures in column F, rows 2-7,
# for both the "展览" and "全部类型" sheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2112
    3 = 613
    4 = 1490
    5 = 7102
    6 = 175
    7 = 131
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
